# Fixed bug in DCOPF to consider storage
#
# 1) "bus" sheet: update the diagram x/y text coordinates for Bus 1/2/3.
#    These are stored as TEXT (not numbers) in the workbook, so we force
#    text formatting before assigning, then clear the formatting again so
#    the cell keeps its original (default) style - only its stored text
#    value changes, exactly like the source diff.
# 2) "battery" sheet: rename the battery from "batt" to "batt1@Bus 3".
# 3) "branch" sheet: fill in real X (and one R) reactance values that were
#    previously placeholder near-zero values, so DCOPF now sees storage.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---- bus sheet : x / y diagram coordinates (stored as text) ----
$wsBus = $wb.Worksheets.Item("bus")

Set-TextValue $wsBus.Range("I2") "-1191.0"   # Bus 1 x
Set-TextValue $wsBus.Range("J2") "-716.0"    # Bus 1 y

Set-TextValue $wsBus.Range("I3") "-1457.0"   # Bus 2 x
Set-TextValue $wsBus.Range("J3") "-861.0"    # Bus 2 y

Set-TextValue $wsBus.Range("I4") "-1319.0"   # Bus 3 x
Set-TextValue $wsBus.Range("J4") "-574.0"    # Bus 3 y

# ---- battery sheet : rename "batt" -> "batt1@Bus 3" (text) ----
$wsBattery = $wb.Worksheets.Item("battery")
Set-TextValue $wsBattery.Range("B2") "batt1@Bus 3"

# ---- branch sheet : real X/R reactance values instead of ~0 ----
$wsBranch = $wb.Worksheets.Item("branch")

$wsBranch.Range("J2").Value = 0.05   # Branch 1 (Bus 3 -> Bus 1) X
$wsBranch.Range("J3").Value = 0.08   # Branch 1 (Bus 1 -> Bus 2) X
$wsBranch.Range("I4").Value = 0.01   # Branch 1 (Bus 3 -> Bus 2) R
$wsBranch.Range("J4").Value = 0.06   # Branch 1 (Bus 3 -> Bus 2) X
